$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Passed"

$ws.Range("C2:C3").Select()
